$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A, B, C) to match the new narrower layout.
# ColumnWidth values chosen so the stored OOXML "width" lands as close as
# possible to the target widths (24 / 13.46484375 / 4.73046875), given the
# runtime's internal pixel-snapping of ColumnWidth.
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 3.8333333333333335

# Update the "Mean flow rate" values in column B (rows 2-9) with the
# recalculated figures.
$ws.Range("B2").Value = 54.951144613083095
$ws.Range("B3").Value = 26.309799551422486
$ws.Range("B4").Value = 18.404036144300758
$ws.Range("B5").Value = 45.598273085931588
$ws.Range("B6").Value = 35.846945556233649
$ws.Range("B7").Value = 12.194011847007653
$ws.Range("B8").Value = 40.594365598338776
$ws.Range("B9").Value = 16.556238448663667
